$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column headers (G: url, H: progress)
$ws.Range("G1").Value = "url"
$ws.Range("H1").Value = "progress"

# New data rows for columns G (string) and H (numeric progress)
$ws.Range("G2").Value = "consent"
$ws.Range("H2").Value = 1

$ws.Range("G3").Value = "family"
$ws.Range("H3").Value = 2

$ws.Range("G4").Value = "breast"
$ws.Range("H4").Value = 3

$ws.Range("G5").Value = "ovarian"
$ws.Range("H5").Value = 3

$ws.Range("G6").Value = "grandmother"
$ws.Range("H6").Value = 4

$ws.Range("G7").Value = "aunt"
$ws.Range("H7").Value = 4

$ws.Range("G8").Value = "niece"
$ws.Range("H8").Value = 4

$ws.Range("G9").Value = "halfsister"
$ws.Range("H9").Value = 4

$ws.Range("G10").Value = "history"
$ws.Range("H10").Value = 5

# Update the sheet view: drop the frozen/scrolled top-left cell and move
# the selection to the cell right after the newly-populated data (H11).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H11").Select()
